# Remove the "Costsheet" (row 23) and "Timesheet" (row 24) rule rows from the
# Drools form business-process rules table on Sheet1. Deleting the two whole
# rows shifts the "ArrestWarrant" row (and the blank rows below it) up by two,
# and the sheet's used range shrinks from A1:K28 to A1:K26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remember the hyperlink that lives in column I of the "ArrestWarrant" row
# (currently row 25) so we can recreate it at its new location after the
# row shift - EntireRow.Delete() does not relocate Hyperlink objects itself.
$hyperlinksBefore = @()
foreach ($hl in $ws.Hyperlinks) {
    $hyperlinksBefore += $hl
}

$arrestWarrantLink = $null
foreach ($hl in $hyperlinksBefore) {
    if ($hl.Range.Address() -eq "$I$25") {
        $arrestWarrantLink = $hl
    }
}
if ($arrestWarrantLink -eq $null) {
    $arrestWarrantLink = $hyperlinksBefore[1]
}

$linkAddress = $arrestWarrantLink.Address
$linkTextToDisplay = $arrestWarrantLink.TextToDisplay

# Delete the two whole rows (Costsheet = row 23, Timesheet = row 24); rows
# below shift up automatically.
$ws.Range("A23:K24").EntireRow.Delete()

# The stale hyperlink object still references the old $I$25 address (now
# blank); drop it and add a fresh one on the relocated ArrestWarrant row
# (now row 23).
$arrestWarrantLink.Delete()
$ws.Hyperlinks.Add($ws.Range("I23"), $linkAddress, "", "", $linkTextToDisplay)

# Re-adding the hyperlink resets I23's style to a generic "Hyperlink" look;
# restore the original formatting (shared with I20/I22) via a formats-only
# paste so no stray style entries linger on the cell.
$ws.Range("I20").Copy()
$ws.Range("I23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Put the selection where the saved file shows it afterwards.
$ws.Range("C23").Select() | Out-Null
